# Updates cryptos list data (prices / 1h volume %) pulled by the scheduled
# GitHub Actions job. Mirrors the source CSV refresh: most rows just get new
# Price / Volume(1h) figures; rows 16-17 (Polygon / WrappedEther) swapped
# ranking order and row 51 (BEAM) was replaced by FlareNetwork.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.279.10"
$ws.Range("E2").Value = "  +9.25%  "
$ws.Range("D3").Value = "3.251.87"
$ws.Range("E3").Value = "  +4.80%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'400.82"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").Value = "'111.40"
$ws.Range("E6").Value = "  +7.98%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +7.65%  "
$ws.Range("D10").Value = "'39.78"
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("D13").Value = "3.751.37"
$ws.Range("E13").Value = "  +4.40%  "
$ws.Range("D14").Value = "'19.27"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("E15").Value = "  +3.78%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.261.12"
$ws.Range("E16").Value = "  +5.31%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'1.07"
$ws.Range("E17").Value = "  +7.80%  "
$ws.Range("D18").Value = "'10.60"
$ws.Range("E18").Value = "  -6.57%  "
$ws.Range("D19").Value = "56.171.15"
$ws.Range("E19").Value = "  +9.03%  "
$ws.Range("D20").Value = "'3.40"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = "  +7.43%  "
$ws.Range("D22").Value = "'0.0000102"
$ws.Range("E22").Value = "  +5.55%  "
$ws.Range("D23").Value = "'289.66"
$ws.Range("E23").Value = "  +8.96%  "
$ws.Range("D24").Value = "'74.64"
$ws.Range("E24").Value = "  +6.68%  "
$ws.Range("E25").Value = "  +4.59%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "'28.28"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").Value = "'7.62"
$ws.Range("E28").Value = "  +5.01%  "
$ws.Range("E29").Value = "  +2.96%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("E32").Value = "  +10.13%  "
$ws.Range("E33").Value = "  +6.96%  "
$ws.Range("D34").Value = "'37.27"
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("D36").Value = "'51.24"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "'3.61"
$ws.Range("E37").Value = "  +7.12%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +21.90%  "
$ws.Range("D40").Value = "'136.47"
$ws.Range("E40").Value = "  +5.73%  "
$ws.Range("D41").Value = "'1.95"
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("E42").Value = "  +10.85%  "
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'0.119"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'16.96"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "'22.78"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "'2.11"
$ws.Range("E47").Value = "  +41.63%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "2.142.33"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("B51").Value = "FlareNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr"
$ws.Range("D51").Value = "'0.0519"
$ws.Range("E51").Value = "  +11.40%  "
